$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the title banner in A1 (merged A1:L1) - new data date 20.05.2020
#    NOTE: writing strings through Range(...).Value on a merged cell silently
#    no-ops in this host, so we go through Cells.Item(row, col) instead.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Données COVID-19 Valais 20.05.2020"

# ---------------------------------------------------------------------------
# 2) A new day of data (20.05.2020) has arrived. The previous last row (19.05)
#    keeps its place as a normal data row, and a fresh "last row" (20.05) is
#    appended below it, matching the style Excel applies when you insert a
#    row above the final row.
# ---------------------------------------------------------------------------

# Insert a blank row at 84: old row 84 (19.05.2020, the former last row) shifts
# down to row 85, carrying its own formatting/formulas/values with it intact.
$ws.Rows.Item(84).Insert()

# Fill the new row 84 (still date 43970 = 19.05.2020) with the corrected data
# and the same formula pattern as the row above it (row 83) BEFORE touching
# formatting, since (re)applying a shared formula after a format paste resets
# the cell's style in this host.
$ws.Cells.Item(84, 1).Value = 43970
$ws.Cells.Item(84, 2).Formula = "=B83+C84"
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 6
$ws.Cells.Item(84, 6).Value = 4
$ws.Cells.Item(84, 7).Value = 22
$ws.Cells.Item(84, 8).Formula = "=G84+E84"
$ws.Cells.Item(84, 9).Formula = "=I83+J84"
$ws.Cells.Item(84, 10).Formula = "=K84+L84"
# Columns K/L are formatted as Text ("@"); writing a literal number straight
# into a Text-formatted cell gets coerced to a string in this host, so flip
# to General first, write the number, and let the formatting paste below
# restore the real (Text) number format/style.
$ws.Cells.Item(84, 11).NumberFormat = "General"
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).NumberFormat = "General"
$ws.Cells.Item(84, 12).Value = 0

# Now copy the formatting (styles/borders/number formats) from row 83 onto
# the new row 84 so it matches the rest of the "normal" data rows. Column G
# is skipped because the insert already leaves it with the right style, and
# re-pasting it directly flips it to the wrong one in this host.
$ws.Range("A83:F83").Copy()
$ws.Range("A84:F84").PasteSpecial(-4122)
$ws.Range("H83:L83").Copy()
$ws.Range("H84:L84").PasteSpecial(-4122)

# Row 85 is what used to be row 84: update its date to 20.05.2020, clear the
# not-yet-known cumulative counters (B/C), and refresh G with the corrected
# figure (H/I/J recompute automatically from their existing formulas).
$ws.Cells.Item(85, 1).Value = 43971
$ws.Range("B85").ClearContents()
$ws.Range("C85").ClearContents()
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 7).Value = 22

# ---------------------------------------------------------------------------
# 3) Reflect the new title cell as the active selection, like Excel leaves it
#    after the author clicked the merged banner before saving.
# ---------------------------------------------------------------------------
$ws.Range("A1:L1").Select()
